# Correcciones lógicas y gráficas
# - Actualiza el historial de pedidos de Alberto Hurtado (agrega el pedido "7").
# - Agrega un nuevo cliente (Javiera Cabrera) con su historial de pedidos ("6").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Alberto Hurtado ahora también tiene el pedido "7" en su historial.
$ws.Range("D2").Value = "0,7,"

# Nueva fila para la clienta Javiera Cabrera.
$ws.Range("A7").Value = "Javiera Cabrera"
$ws.Range("B7").Value = "javieracabrera14@gmail.com"
$ws.Range("C7").Value = "9-66666666"

# El valor "6," se interpretaría como el número 6 si se escribe directamente
# (Excel recorta la coma final de un literal puramente numérico). Se construye
# como texto vía fórmula y se pega como valor para preservar la coma literal
# sin tocar el formato de la celda (queda con el estilo general, como el resto
# de la tabla).
$ws.Range("D7").Formula = '="6,"'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D7").PasteSpecial("xlPasteValues") | Out-Null
$excel.CutCopyMode = $false
